# Nieuwe data toegevoegd via Streamlit op 2024-12-03 18:04:40
# Append a new record row (63) to the CompaNanny database sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 63

$ws.Cells.Item($row, 1).Value = "CompaNanny"
$ws.Cells.Item($row, 2).Value = "CompaNanny Benoordenhout KDV"
$ws.Cells.Item($row, 3).Value = "KDV"

# Keep the report date as literal text (matches the other recent rows),
# rather than letting it be auto-converted into a date serial number.
$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "2024-09-23"
$ws.Cells.Item($row, 4).Style = "Normal"

$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 1
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
